$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Rename the second sheet
$ws2.Name = "Supplier and cost"

# Capture the existing "Material / Supplier / Cost per unit" table data
# (currently living at C3:E7) before we move it.
$data = @()
for ($r = 3; $r -le 7; $r++) {
  $row = @()
  for ($c = 3; $c -le 5; $c++) {
    $row += ,($ws2.Cells.Item($r, $c).Value2)
  }
  $data += ,$row
}

# Remove the table definition so clearing/rewriting the underlying
# cells doesn't corrupt its column headers, then clear the old cells.
$lo = $ws2.ListObjects.Item(1)
$lo.Unlist()
$ws2.Range("C3:E7").Clear()

# Write the captured values into the new top-left location A1:C5.
for ($ri = 0; $ri -lt 5; $ri++) {
  for ($ci = 0; $ci -lt 3; $ci++) {
    $ws2.Cells.Item($ri + 1, $ci + 1).Value = $data[$ri][$ci]
  }
}

# Recreate the table over the new range, keeping the same name/style.
$newLo = $ws2.ListObjects.Add(1, $ws2.Range("A1:C5"), 0, 1)
$newLo.Name = "Table2"
$newLo.TableStyle = "TableStyleLight2"

# Resize columns to fit the relocated content.
$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()
$ws2.Columns.Item(3).AutoFit()

# The supplier/cost sheet is now the active one.
$ws2.Activate()
$ws2.Range("E29").Select()
